$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3062484
$ws.Cells.Item(137, 9).Value = 1667872.2
$ws.Cells.Item(137, 10).Value = 5264502.5
$ws.Cells.Item(137, 11).Value = 5003616.6
$ws.Cells.Item(137, 12).Value = 15793507.5
$ws.Cells.Item(137, 13).Value = -5001066.6
$ws.Cells.Item(137, 14).Value = -15798607.5
$ws.Cells.Item(138, 8).Value = 3024.425
$ws.Cells.Item(138, 9).Value = 2466.077
$ws.Cells.Item(138, 10).Value = 3293.2593
$ws.Cells.Item(138, 11).Value = 7398.231000000001
$ws.Cells.Item(138, 12).Value = 9879.777900000001
$ws.Cells.Item(138, 13).Value = -2258.231000000001
$ws.Cells.Item(138, 14).Value = -20159.7779

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20772.38
$ws.Cells.Item(32, 9).Value = 4435.0894
$ws.Cells.Item(32, 11).Value = 4435.0894
$ws.Cells.Item(32, 13).Value = -4148.0894
$ws.Cells.Item(74, 8).Value = 866.7451
$ws.Cells.Item(74, 9).Value = 846.27905
$ws.Cells.Item(74, 11).Value = 846.27905
$ws.Cells.Item(74, 13).Value = 27.72095000000002
$ws.Cells.Item(77, 8).Value = 866.7451
$ws.Cells.Item(77, 9).Value = 846.27905
$ws.Cells.Item(77, 11).Value = 4231.39525
$ws.Cells.Item(77, 13).Value = 136.6047500000004
$ws.Cells.Item(132, 8).Value = 123561.09
$ws.Cells.Item(132, 9).Value = 152945.25
$ws.Cells.Item(132, 10).Value = 2351.5
$ws.Cells.Item(132, 11).Value = 458835.75
$ws.Cells.Item(132, 12).Value = 7054.5
$ws.Cells.Item(132, 13).Value = -456305.75
$ws.Cells.Item(132, 14).Value = -12114.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9805567
$ws.Cells.Item(31, 9).Value = 1457.2
$ws.Cells.Item(31, 10).Value = 23811438
$ws.Cells.Item(31, 11).Value = 1457.2
$ws.Cells.Item(31, 12).Value = 23811438
$ws.Cells.Item(31, 13).Value = -1162.2
$ws.Cells.Item(31, 14).Value = -23812028
$ws.Cells.Item(34, 8).Value = 9805567
$ws.Cells.Item(34, 9).Value = 1457.2
$ws.Cells.Item(34, 10).Value = 23811438
$ws.Cells.Item(34, 11).Value = 1457.2
$ws.Cells.Item(34, 12).Value = 23811438
$ws.Cells.Item(34, 13).Value = -1255.2
$ws.Cells.Item(34, 14).Value = -23811842
$ws.Cells.Item(38, 8).Value = 25237.5
$ws.Cells.Item(38, 9).Value = 1975
$ws.Cells.Item(38, 10).Value = 48500
$ws.Cells.Item(38, 11).Value = 1975
$ws.Cells.Item(38, 12).Value = 48500
$ws.Cells.Item(38, 13).Value = -1598
$ws.Cells.Item(38, 14).Value = -49254
$ws.Cells.Item(46, 8).Value = 25237.5
$ws.Cells.Item(46, 9).Value = 1975
$ws.Cells.Item(46, 10).Value = 48500
$ws.Cells.Item(46, 11).Value = 1975
$ws.Cells.Item(46, 12).Value = 48500
$ws.Cells.Item(46, 13).Value = -1764
$ws.Cells.Item(46, 14).Value = -48922
$ws.Cells.Item(58, 8).Value = 968.97675
$ws.Cells.Item(58, 9).Value = 1030.4375
$ws.Cells.Item(58, 10).Value = 790.1818
$ws.Cells.Item(58, 11).Value = 1030.4375
$ws.Cells.Item(58, 12).Value = 790.1818
$ws.Cells.Item(58, 13).Value = -827.4375
$ws.Cells.Item(58, 14).Value = -1196.1818
$ws.Cells.Item(134, 8).Value = 5411130.5
$ws.Cells.Item(134, 9).Value = 6049.3335
$ws.Cells.Item(134, 11).Value = 18148.0005
$ws.Cells.Item(134, 13).Value = -15613.0005
$ws.Cells.Item(136, 8).Value = 968.97675
$ws.Cells.Item(136, 9).Value = 1030.4375
$ws.Cells.Item(136, 10).Value = 790.1818
$ws.Cells.Item(136, 11).Value = 3091.3125
$ws.Cells.Item(136, 12).Value = 2370.5454
$ws.Cells.Item(136, 13).Value = -541.3125
$ws.Cells.Item(136, 14).Value = -7470.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 66703.47
$ws.Cells.Item(2, 9).Value = 100033.2
$ws.Cells.Item(2, 10).Value = 44
$ws.Cells.Item(2, 11).Value = 600199.2
$ws.Cells.Item(2, 12).Value = 264
$ws.Cells.Item(2, 13).Value = -600086.2
$ws.Cells.Item(2, 14).Value = -490
$ws.Cells.Item(5, 8).Value = 1212866.6
$ws.Cells.Item(5, 9).Value = 820
$ws.Cells.Item(5, 11).Value = 2460
$ws.Cells.Item(5, 13).Value = -2348
$ws.Cells.Item(15, 8).Value = 116.666664
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()
$ws.Cells.Item(21, 8).Value = 135.5
$ws.Cells.Item(21, 9).Value = 135.5
$ws.Cells.Item(21, 11).Value = 406.5
$ws.Cells.Item(21, 13).Value = -233.5
$ws.Cells.Item(33, 8).Value = 548.5833
$ws.Cells.Item(33, 9).Value = 487.7143
$ws.Cells.Item(33, 10).Value = 633.8
$ws.Cells.Item(33, 11).Value = 2926.2858
$ws.Cells.Item(33, 12).Value = 3802.8
$ws.Cells.Item(33, 13).Value = -2643.2858
$ws.Cells.Item(33, 14).Value = -4368.799999999999
$ws.Cells.Item(44, 8).Value = 1897.375
$ws.Cells.Item(44, 10).Value = 1897.375
$ws.Cells.Item(44, 12).Value = 5692.125
$ws.Cells.Item(44, 14).Value = -6488.125
$ws.Cells.Item(132, 8).Value = 941025.4399999999
$ws.Cells.Item(132, 9).Value = 1197055.1
$ws.Cells.Item(132, 10).Value = 2250
$ws.Cells.Item(132, 11).Value = 10773495.9
$ws.Cells.Item(132, 12).Value = 20250
$ws.Cells.Item(132, 13).Value = -10770965.9
$ws.Cells.Item(132, 14).Value = -25310
$ws.Cells.Item(135, 8).Value = 1212866.6
$ws.Cells.Item(135, 9).Value = 820
$ws.Cells.Item(135, 11).Value = 7380
$ws.Cells.Item(135, 13).Value = -4845
$ws.Cells.Item(137, 8).Value = 2562.2222
$ws.Cells.Item(137, 9).Value = 2151.4285
$ws.Cells.Item(137, 10).Value = 4000
$ws.Cells.Item(137, 11).Value = 6454.2855
$ws.Cells.Item(137, 12).Value = 12000
$ws.Cells.Item(137, 13).Value = -1354.2855
$ws.Cells.Item(137, 14).Value = -22200

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3702.3333
$ws.Cells.Item(122, 9).Value = 4204.2
$ws.Cells.Item(122, 10).Value = 3075
$ws.Cells.Item(122, 11).Value = 12612.6
$ws.Cells.Item(122, 12).Value = 9225
$ws.Cells.Item(122, 13).Value = -10162.6
$ws.Cells.Item(122, 14).Value = -14125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 5005
$ws.Cells.Item(9, 9).Value = 673.3333
$ws.Cells.Item(9, 10).Value = 18000
$ws.Cells.Item(9, 11).Value = 673.3333
$ws.Cells.Item(9, 12).Value = 18000
$ws.Cells.Item(9, 13).Value = -449.3333
$ws.Cells.Item(9, 14).Value = -18448
$ws.Cells.Item(30, 8).Value = 38816
$ws.Cells.Item(30, 9).Value = 38816
$ws.Cells.Item(30, 11).Value = 38816
$ws.Cells.Item(30, 13).Value = -38708
$ws.Cells.Item(35, 8).Value = 3566.6667
$ws.Cells.Item(35, 9).Value = 1850
$ws.Cells.Item(35, 10).Value = 7000
$ws.Cells.Item(35, 11).Value = 1850
$ws.Cells.Item(35, 12).Value = 7000
$ws.Cells.Item(35, 13).Value = -1514
$ws.Cells.Item(35, 14).Value = -7672
$ws.Cells.Item(134, 8).Value = 32493.334
$ws.Cells.Item(134, 10).Value = 32493.334
$ws.Cells.Item(134, 12).Value = 32493.334
$ws.Cells.Item(134, 14).Value = -42633.334
$ws.Cells.Item(138, 8).Value = 42437.43
$ws.Cells.Item(138, 10).Value = 42437.43
$ws.Cells.Item(138, 12).Value = 42437.43
$ws.Cells.Item(138, 14).Value = -52717.43
$ws.Cells.Item(141, 8).Value = 74430
$ws.Cells.Item(141, 10).Value = 74430
$ws.Cells.Item(141, 12).Value = 74430
$ws.Cells.Item(141, 14).Value = -84790

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2882.7715
$ws.Cells.Item(132, 9).Value = 3045.926
$ws.Cells.Item(132, 10).Value = 2332.125
$ws.Cells.Item(132, 11).Value = 9137.778
$ws.Cells.Item(132, 12).Value = 6996.375
$ws.Cells.Item(132, 13).Value = -6607.778
$ws.Cells.Item(132, 14).Value = -12056.375
$ws.Cells.Item(136, 8).Value = 1275.6735
$ws.Cells.Item(136, 9).Value = 1267.6279
$ws.Cells.Item(136, 10).Value = 1333.3334
$ws.Cells.Item(136, 11).Value = 3802.8837
$ws.Cells.Item(136, 12).Value = 4000.0002
$ws.Cells.Item(136, 13).Value = -1252.8837
$ws.Cells.Item(136, 14).Value = -9100.0002

Write-Output "Applied all changes"